# Correct hourly gas shape
# Add two new dwelling type lookup rows (TERRACED) to Sheet1 of the
# lookup_dwelling_type workbook, extending the table from A1:B7 to A1:B9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "TERRACED"

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "TERRACED"

# Match the author's final cursor/selection position recorded in the sheet.
$ws.Range("D9").Select()
